$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shrink the sheet from 27 data rows to 25 (rows 26 and 27 are removed).
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()

# ---------------------------------------------------------------------------
# 2. Row 2 - CENCOSUD / YOGURT GRIEGO
#    A2 is a genuine number (50001179), not a shared string.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 50001179
$ws.Range("B2").Value = "CENCOSUD"
$ws.Range("C2").Value = "484329005"
$ws.Range("D2").Value = "YOGURT GRIEGO LAIVE ORIGINAL ENDULZADO"
$ws.Range("E2").Value = 0.5

# ---------------------------------------------------------------------------
# 3. Rows 3-6 - Daro products.
#    Column A holds numeric-looking codes that must stay TEXT (shared
#    string) even though the column's number format is General, so the
#    apostrophe-prefix trick is used and then the quote-prefix style is
#    replaced by copying the plain, unmodified style from A10 (still style
#    index 6 in the source file, never touched by this script).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "'80000156"
$ws.Range("A4").Value = "'80000157"
$ws.Range("A5").Value = "'80000158"
$ws.Range("A6").Value = "'80000155"
$ws.Range("A10").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B3").Value = "Daro"
$ws.Range("B4").Value = "Daro"
$ws.Range("B5").Value = "Daro"
$ws.Range("B6").Value = "Daro"

$ws.Range("C3").Value = "80000156       "
$ws.Range("C4").Value = "80000157       "
$ws.Range("C5").Value = "80000158       "
$ws.Range("C6").Value = "80000155       "

$ws.Range("D3").Value = "BISCOLATA VENITRIO CREM AVELL 20G/24U/6D"
$ws.Range("D4").Value = "WINERGY BARRA CHOCO MANI 18G/24U/12D"
$ws.Range("D5").Value = "GALL BISCOLATA MOOD RELL CHOC 25G/24U/6D"
$ws.Range("D6").Value = "BISCOLATA CARAM CUBIET CHOCO 20G/24U/12D"

$ws.Range("E3").Value = 0.02
$ws.Range("E4").Value = 0.018
$ws.Range("E5").Value = 0.025
$ws.Range("E6").Value = 0.02

# Copy the body row formatting (row 6, untouched previously) onto the
# D column cells of rows 3-5 so they pick up the same style D6 already had
# instead of the old style inherited from the legacy row 2/3 layout.
$ws.Range("D6").Copy()
$ws.Range("D2:D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the values on top (PasteSpecial of formats only should not have
# touched them, but this keeps value + format both correct regardless of
# paste semantics).
$ws.Range("D2").Value = "YOGURT GRIEGO LAIVE ORIGINAL ENDULZADO"
$ws.Range("D3").Value = "BISCOLATA VENITRIO CREM AVELL 20G/24U/6D"
$ws.Range("D4").Value = "WINERGY BARRA CHOCO MANI 18G/24U/12D"
$ws.Range("D5").Value = "GALL BISCOLATA MOOD RELL CHOC 25G/24U/6D"
$ws.Range("D6").Value = "BISCOLATA CARAM CUBIET CHOCO 20G/24U/12D"

# ---------------------------------------------------------------------------
# 4. Row 7 - DIJISA / MEZCLA LACTEA NUTRILAC
#    A7 is a genuine number (50001045).
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 50001045
$ws.Range("B7").Value = "DIJISA"
$ws.Range("C7").Value = "047329              "
$ws.Range("D7").Value = "MEZCLA LACTEA NUTRILAC CJ*480GR"
$ws.Range("E7").Value = 0.48

# ---------------------------------------------------------------------------
# 5. Selection bookkeeping to mirror the saved view state.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Select()
